$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to keep its text (string) representation even when the
    # supplied value looks like a plain number (e.g. "212.83"), then restore
    # the cell's original (default) style so no stray formatting is left behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.262.30"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.594.28"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "212.83"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.58%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.245"
$ws.Range("E8").Value = "  -0.29%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.55%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -2.54%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0852"
$ws.Range("E11").Value = "  +0.60%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.819.70"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.597.22"
$ws.Range("E13").Value = "  +0.50%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -2.58%  "

# Row 16 - Litecoin
Set-TextValue "D16" "63.95"
$ws.Range("E16").Value = "  -1.19%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "26.258.60"

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.64%  "

# Row 19 - now Chainlink (was BitcoinCash)
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "7.36"
$ws.Range("E19").Value = "  -1.74%  "

# Row 20 - now BitcoinCash (was Chainlink)
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "213.82"
$ws.Range("E20").Value = "  +0.91%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.12%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -0.07%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.05"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -3.00%  "

# Row 25 - Monero
Set-TextValue "D25" "144.99"
$ws.Range("E25").Value = "  +0.12%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.03%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -1.62%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.59%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -2.57%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.47%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -0.55%  "

# Row 33 - Maker
Set-TextValue "D33" "1.418.03"
$ws.Range("E33").Value = "  +5.80%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "2.99"
$ws.Range("E34").Value = "  +0.12%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.52%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  -1.03%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "0.580"
$ws.Range("E37").Value = "  -3.08%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.97%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  +0.66%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - WEMIXToken
Set-TextValue "D42" "0.966"
$ws.Range("E42").Value = "  -9.60%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +1.02%  "

# Row 44 - TrustWalletToken
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.730.85"

# Row 46 - Aave
Set-TextValue "D46" "60.95"
$ws.Range("E46").Value = "  -1.08%  "

# Row 47 - Quant
Set-TextValue "D47" "86.89"
$ws.Range("E47").Value = "  -1.26%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.26%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.58%  "

# Row 50 - Algorand
Set-TextValue "D50" "0.0956"
$ws.Range("E50").Value = "  -2.83%  "

# Row 51 - USDD
Set-TextValue "D51" "0.999"
$ws.Range("E51").Value = "  -0.03%  "
